# Correlation Analysis (Not Done yet)
# Applies date corrections, renames "Chris" -> "Christopher", and updates
# the sheet view's selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct experiment dates in column G ---
$ws.Range("G3").Value = 20170430

$ws.Range("G5").Value = 20170503
$ws.Range("G6").Value = 20170503
$ws.Range("G7").Value = 20170503
$ws.Range("G8").Value = 20170503

$ws.Range("G9").Value = 20170509
$ws.Range("G10").Value = 20170509
$ws.Range("G11").Value = 20170509
$ws.Range("G12").Value = 20170509

$ws.Range("G13").Value = 20170517
$ws.Range("G14").Value = 20170517
$ws.Range("G15").Value = 20170517
$ws.Range("G16").Value = 20170517

$ws.Range("G17").Value = 20170519
$ws.Range("G18").Value = 20170519

$ws.Range("G19").Value = 20170525

# --- Fix participant first name ---
$ws.Range("B19").Value = "Christopher"

# --- Update selection / scroll position ---
$ws.Range("A19").Select()
